# Apply monthly data corrections to the 沪铝期货价格 (Shanghai Aluminum Futures Price) sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 21023
$ws.Range("B3").Value = 21048
$ws.Range("C3").Value = 20969
$ws.Range("C4").Value = 20784
$ws.Range("C6").Value = 20612
$ws.Range("C11").Value = 20613
$ws.Range("C12").Value = 20121
$ws.Range("C14").Value = 20782
$ws.Range("C17").Value = 19404
$ws.Range("C19").Value = 20678
$ws.Range("C20").Value = 20904
$ws.Range("C22").Value = 19271
$ws.Range("C23").Value = 18841
